# Reorder the two leave-request rows on Sheet4 (swap rows 4 and 5 - values
# AND formatting move together, as a genuine row swap), then make Sheet4 the
# active sheet with its whole 5th row selected (Sheet3 was previously the
# selected/active tab, with H7 selected).

$wb  = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("Sheet4")

# --- Swap the contents+formatting of row 4 and row 5 on Sheet4 -------------
# Range.Copy(Destination) carries the source formatting (number format,
# font, quote-prefix, ...) to the destination, unlike a plain .Value copy.
# Use a scratch row far below the data as temporary storage so the two
# source rows don't clobber one another, then remove that scratch row with
# Delete (not just ClearContents) so it doesn't linger in the sheet's used
# range / dimension.
$ws4.Range("A4:D4").Copy($ws4.Range("A20:D20"))
$ws4.Range("A5:D5").Copy($ws4.Range("A4:D4"))
$ws4.Range("A20:D20").Copy($ws4.Range("A5:D5"))
$ws4.Rows("20:20").Delete()

# --- Update the active sheet / selection ------------------------------------
# Sheet4 becomes the selected tab, with its whole 5th row as the active
# selection (A5:XFD5).
$ws4.Activate() | Out-Null
$ws4.Rows("5:5").Select() | Out-Null
